{"js": "// The three title-block paragraphs (Title, Author, Abstract) were each\n// split word-by-word across many runs (with separate single-space runs in\n// between). Collapse each paragraph's text into a single run by replacing\n// the paragraph's full range with its own (unchanged) text content.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\nconst targetStyles = new Set([\"Title\", \"Author\", \"Abstract\"]);\n\nfor (const paragraph of paragraphs.items) {\n  if (targetStyles.has(paragraph.style)) {\n    const fullText = paragraph.text;\n    paragraph.getRange().insertText(fullText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# The three title-block paragraphs (Title, Author, Abstract) were each\n# split word-by-word across many runs (with separate single-space runs in\n# between). Collapse each paragraph's text into a single run by running a\n# Find/Replace of the paragraph's own text over its own range - this makes\n# Word rewrite the range as one run instead of leaving the old word-by-word\n# run split in place.\n$d = $word.ActiveDocument\n$targetStyles = @(\"Title\", \"Author\", \"Abstract\")\n\nforeach ($p in $d.Paragraphs) {\n    if ($targetStyles -contains $p.Style.NameLocal) {\n        $r = $p.Range\n        $plainText = $r.Text.TrimEnd(\"`r\")\n        [void]$r.Find.Execute($plainText, $false, $false, $false, $false, $false, $true, 1, $false, $plainText, 2)\n    }\n}\n"}
